$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 187 - this shifts rows 187:246 down to 188:247,
# inheriting the formatting of the row above (row 186), matching Excel's
# default row-insert behavior.
$ws.Rows.Item(187).Insert()

# Copy the common/static columns from the row above (row 186) into the new row 187
$ws.Range("A187").Value = $ws.Range("A186").Value2
$ws.Range("B187").Value = $ws.Range("B186").Value2
$ws.Range("C187").Value = $ws.Range("C186").Value2
$ws.Range("E187").Value = $ws.Range("E186").Value2
$ws.Range("F187").Value = $ws.Range("F186").Value2
$ws.Range("G187").Value = $ws.Range("G186").Value2
$ws.Range("H187").Value = $ws.Range("H186").Value2
$ws.Range("I187").Value = $ws.Range("I186").Value2
$ws.Range("N187").Value = $ws.Range("N186").Value2
$ws.Range("O187").Value = $ws.Range("O186").Value2
$ws.Range("Q187").Value = $ws.Range("Q186").Value2
$ws.Range("R187").Value = $ws.Range("R186").Value2

# Set the new row's own data values
$ws.Range("D187").Value = 45146
$ws.Range("J187").Value = 2000
$ws.Range("K187").Value = 2500
$ws.Range("L187").Value = 3000
$ws.Range("M187").Value = 2750
$ws.Range("P187").Value = 1833

$ws.Range("A1").Select()
